$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the time_taken (F) column timestamps on the "data" sheet ---
$newTimes = @(
    "2021-10-05 14:33:26.691116",
    "2021-10-05 14:33:26.691122",
    "2021-10-05 14:33:26.691125",
    "2021-10-05 14:33:26.691127",
    "2021-10-05 14:33:26.691130",
    "2021-10-05 14:33:26.691132",
    "2021-10-05 14:33:26.691134",
    "2021-10-05 14:33:26.691136",
    "2021-10-05 14:33:26.691138",
    "2021-10-05 14:33:26.691140",
    "2021-10-05 14:33:26.691142",
    "2021-10-05 14:33:26.691144",
    "2021-10-05 14:33:26.691145",
    "2021-10-05 14:33:26.691147",
    "2021-10-05 14:33:26.691149",
    "2021-10-05 14:33:26.691151",
    "2021-10-05 14:33:26.691153",
    "2021-10-05 14:33:26.691155",
    "2021-10-05 14:33:26.691157",
    "2021-10-05 14:33:26.691159",
    "2021-10-05 14:33:26.691161",
    "2021-10-05 14:33:26.691164",
    "2021-10-05 14:33:26.691165",
    "2021-10-05 14:33:26.691167",
    "2021-10-05 14:33:26.691170",
    "2021-10-05 14:33:26.691172",
    "2021-10-05 14:33:26.691174",
    "2021-10-05 14:33:26.691176",
    "2021-10-05 14:33:26.691178",
    "2021-10-05 14:33:26.691180",
    "2021-10-05 14:33:26.691181",
    "2021-10-05 14:33:26.691183",
    "2021-10-05 14:33:26.691185",
    "2021-10-05 14:33:26.691188",
    "2021-10-05 14:33:26.691190",
    "2021-10-05 14:33:26.691192",
    "2021-10-05 14:33:26.691194",
    "2021-10-05 14:33:26.691196",
    "2021-10-05 14:33:26.691198",
    "2021-10-05 14:33:26.691200",
    "2021-10-05 14:33:26.691202",
    "2021-10-05 14:33:26.691204",
    "2021-10-05 14:33:26.691206",
    "2021-10-05 14:33:26.691215",
    "2021-10-05 14:33:26.691217",
    "2021-10-05 14:33:26.691219",
    "2021-10-05 14:33:26.691221",
    "2021-10-05 14:33:26.691223",
    "2021-10-05 14:33:26.691225",
    "2021-10-05 14:33:26.691227",
    "2021-10-05 14:33:26.691229",
    "2021-10-05 14:33:26.691231",
    "2021-10-05 14:33:26.691234",
    "2021-10-05 14:33:26.691236",
    "2021-10-05 14:33:26.691238",
    "2021-10-05 14:33:26.691239",
    "2021-10-05 14:33:26.691242",
    "2021-10-05 14:33:26.691245",
    "2021-10-05 14:33:26.691248",
    "2021-10-05 14:33:26.691252",
    "2021-10-05 14:33:26.691255",
    "2021-10-05 14:33:26.691258",
    "2021-10-05 14:33:26.691261",
    "2021-10-05 14:33:26.691263",
    "2021-10-05 14:33:26.691266",
    "2021-10-05 14:33:26.691268",
    "2021-10-05 14:33:26.691271",
    "2021-10-05 14:33:26.691273",
    "2021-10-05 14:33:26.691275",
    "2021-10-05 14:33:26.691277",
    "2021-10-05 14:33:26.691279",
    "2021-10-05 14:33:26.691281",
    "2021-10-05 14:33:26.691283",
    "2021-10-05 14:33:26.691285",
    "2021-10-05 14:33:26.691287",
    "2021-10-05 14:33:26.691289",
    "2021-10-05 14:33:26.691293",
    "2021-10-05 14:33:26.691295",
    "2021-10-05 14:33:26.691297",
    "2021-10-05 14:33:26.691299",
    "2021-10-05 14:33:26.691302",
    "2021-10-05 14:33:26.691304",
    "2021-10-05 14:33:26.691306",
    "2021-10-05 14:33:26.691308",
    "2021-10-05 14:33:26.691310",
    "2021-10-05 14:33:26.691312",
    "2021-10-05 14:33:26.691314",
    "2021-10-05 14:33:26.691316",
    "2021-10-05 14:33:26.691318",
    "2021-10-05 14:33:26.691320",
    "2021-10-05 14:33:26.691322",
    "2021-10-05 14:33:26.691325"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" worksheet, placed right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$metaSheet.Name = "metadata"

# Header row (B1:G1) -- values
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row (A2:G2) -- values
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Cholestasis"
$metaSheet.Cells.Item(2, 3).Value = 78
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "0.204"
$metaSheet.Cells.Item(2, 4).Style = "Normal"
$metaSheet.Cells.Item(2, 5).Value = "2021-09-27T07:38:56.058169Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:33:26.688556"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/78/?format=json"

# --- Formatting: reuse the same header/index styling already used on "data" ---
$ws.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the originally-active "data" sheet as the active tab (adding a
# worksheet normally activates it, which we don't want here).
$ws.Activate()

Write-Output "done"
